# Final & last changes #SNET-26
#
# The original workbook has a single sheet "Personas" that contains both a
# trimmed set of columns (A:J) used for the public-facing form responses and
# two extra columns (K:L, "Código de promoción" / "Otra") used internally.
#
# The edit:
#   1. Renames "Personas" -> "Propuesta" (the trimmed, public view).
#   2. Duplicates it (before trimming) into a new sheet "Original" that keeps
#      every column (A:L) exactly as before.
#   3. Removes the no-longer-needed columns K:L from "Propuesta".
#   4. Leaves "Propuesta" as the active/selected tab (cell D17 selected) and
#      "Original" selected at cell A42.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the existing sheet -----------------------------------------
$propuesta = $wb.ActiveSheet
$propuesta.Name = "Propuesta"

# --- 2. Duplicate it (full A:L content) as "Original" ----------------------
$propuesta.Copy([System.Reflection.Missing]::Value, $propuesta) | Out-Null
$original = $wb.Worksheets.Item(2)
$original.Name = "Original"

# --- 3. Drop the promo-code / other columns from "Propuesta" --------------
$propuesta.Range("K1:L25").EntireColumn.Delete() | Out-Null

# --- 4. Restore tab/selection state ----------------------------------------
$original.Activate()
$original.Range("A42").Select() | Out-Null

$propuesta.Activate()
$propuesta.Range("D17").Select() | Out-Null
